$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the first paragraph
#    ("The objective of this lab is to give you practice working
#    with: ") down to the start of the "The assignment operator"
#    bullet paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$findRange = $d.Content
$findRange.Find.Execute("The assignment operato")
$targetPara = $findRange.Paragraphs(1)
$target = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $target)

# ---------------------------------------------------------------------
# 2. Fix the header: "CS 133JS" -> "CS133JS"
# ---------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("CS 133JS", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "CS133JS", 2)
        }
    }
}
